$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 373
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 8

$ws.Range("B5").Value = 0.9325
$ws.Range("C5").Value = 0.0475
$ws.Range("D5").Value = 0.02
